# Refresh the "cryptos" price/volume table (GitHub Actions bot update).
#
# Note: several Price (column D) values look like plain numbers (e.g.
# "32.27") and Excel's COM layer auto-converts such strings to real
# numbers on assignment, which would change the cell's stored type away
# from the original plain-text cell. To keep those cells as text (as in
# the source file), we set NumberFormat to "@" (Text) immediately before
# assigning the value, then reset the style back to "Normal" afterwards
# -- this preserves the text type without leaving a lingering custom
# number-format style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.762.53'
$ws.Range('E2').Value = '  +4.21%  '

# Row 3
$ws.Range('D3').Value = '2.266.53'
$ws.Range('E3').Value = '  +2.30%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.07%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.27'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.48%  '

# Row 7
$ws.Range('E7').Value = '  +3.44%  '

# Row 8
$ws.Range('E8').Value = '  -0.06%  '

# Row 9
$ws.Range('E9').Value = '  +2.40%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.10%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.44'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.04%  '

# Row 12
$ws.Range('E12').Value = '  +1.72%  '

# Row 13
$ws.Range('E13').Value = '  +1.23%  '

# Row 14
$ws.Range('E14').Value = '  +3.12%  '

# Row 15
$ws.Range('D15').Value = '2.617.39'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.26%  '

# Row 17
$ws.Range('D17').Value = '2.266.79'
$ws.Range('E17').Value = '  +2.42%  '

# Row 18
$ws.Range('E18').Value = '  +3.46%  '

# Row 19
$ws.Range('D19').Value = '41.694.79'
$ws.Range('E19').Value = '  +4.21%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.97%  '

# Row 21
$ws.Range('D21').Value = '0.0₃0901'
$ws.Range('E21').Value = '  +1.66%  '

# Row 22
$ws.Range('E22').Value = '  +2.35%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.59%  '

# Row 25
$ws.Range('E25').Value = '  +4.07%  '

# Row 26
$ws.Range('E26').Value = '  +0.05%  '

# Row 27
$ws.Range('E27').Value = '  +5.04%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.82%  '

# Row 29
$ws.Range('E29').Value = '  +2.14%  '

# Row 30
$ws.Range('E30').Value = '  +2.12%  '

# Row 31
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.68%  '

# Row 32
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '159.97'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.65%  '

# Row 33
$ws.Range('E33').Value = '  -0.02%  '

# Row 34
$ws.Range('E34').Value = '  +4.07%  '

# Row 35
$ws.Range('E35').Value = '  +4.23%  '

# Row 36
$ws.Range('E36').Value = '  -1.55%  '

# Row 37
$ws.Range('E37').Value = '  +1.92%  '

# Row 38
$ws.Range('E38').Value = '  +2.67%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.54%  '

# Row 40
$ws.Range('E40').Value = '  +3.60%  '

# Row 41
$ws.Range('E41').Value = '  +2.51%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.24%  '

# Row 43
$ws.Range('D43').Value = '2.059.49'
$ws.Range('E43').Value = '  -0.78%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.32%  '

# Row 45
$ws.Range('E45').Value = '  +2.62%  '

# Row 46
$ws.Range('E46').Value = '  +2.24%  '

# Row 47
$ws.Range('E47').Value = '  +3.50%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.04'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.81%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '73.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.38%  '

# Row 50
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.81%  '

# Row 51
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.15'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.80%  '
